$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cells for the "fixed boundary condition" section (columns G/H)
# String values are set in an order that reproduces the shared-string table
# order of the target workbook: l²/EI , E, I, L
$ws.Range("G2").Value = "l²/EI "

$ws.Range("G3").Formula = "=H5^2/H6/H7"

$ws.Range("G6").Value = "E"
$ws.Range("H6").Value = 1

$ws.Range("G7").Value = "I"
$ws.Range("H7").Value = 1

$ws.Range("G5").Value = "L"
$ws.Range("H5").Value = 1

# Update the selected cell to match the new active cell
$ws.Range("H8").Select()
